$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as text (matches inline-string cells in source)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "25.736.58"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.618.65"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "214.57"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "0.5080"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.2564"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "0.06355"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "0.07771"
$ws.Range("D12").Value = "4.238"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "1.622.82"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "1.838.44"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "0.5545"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "63.47"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "0.0₅7497"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "25.722.53"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("D20").Value = "193.57"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "4.341"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "9.743"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "5.957"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D26").Value = "140.69"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "0.1262"
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").Value = "6.717"
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("D29").Value = "15.39"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "1.232"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "0.04849"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "3.288"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "3.171"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "1.546"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "2.364"
$ws.Range("D36").Value = "0.8914"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "1.124.20"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.528"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.5477"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").Value = "0.01556"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "0.7906"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "96.94"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "1.762.05"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("D47").Value = "0.4409"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "54.56"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "0.05063"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("D50").Value = "7.515"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "0.9983"
$ws.Range("E51").Value = "  -0.27%  "
